# Update gh-pages to output generated at 456a3b4
# Refresh "想去人数" (F) and, where sold-out, "最低票价" (G) counts across sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value  = 267
$ws.Range("F4").Value  = 1083
$ws.Range("F5").Value  = 2522
$ws.Range("F8").Value  = 52
$ws.Range("F9").Value  = 233
$ws.Range("F12").Value = 82
$ws.Range("F13").Value = 106
$ws.Range("F14").Value = 1461

# ---- Sheet "演出" ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F10").Value = 16
$ws.Range("F19").Value = 49

# ---- Sheet "本地生活" ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6342
$ws.Range("G2").Value = 0
$ws.Range("F4").Value = 2009
$ws.Range("F5").Value = 236

# ---- Sheet "全部类型" ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 6342
$ws.Range("G2").Value  = 0
$ws.Range("F4").Value  = 2009
$ws.Range("F5").Value  = 236
$ws.Range("F11").Value = 267
$ws.Range("F12").Value = 1083
$ws.Range("F16").Value = 2522
$ws.Range("F19").Value = 16
$ws.Range("F23").Value = 52
$ws.Range("F24").Value = 233
$ws.Range("F28").Value = 82
$ws.Range("F29").Value = 106
$ws.Range("F31").Value = 1461
$ws.Range("F39").Value = 49
